$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Scanner -> Histology_Lab)
$ws.Name = "Histology_Lab"

# Update headers (row 1): Number/Content/Location/Log Date/Log Time/Type
# -> Student ID/Location/Log Date/Log Time/Number
$ws.Range("A1").Value = "Student ID"
$ws.Range("B1").Value = "Location"
$ws.Range("C1").Value = "Log Date"
$ws.Range("D1").Value = "Log Time"
$ws.Range("E1").Value = "Number"

# Update the data row. Some of these values look numeric/date-like to
# Excel's auto-detection (the student id, the date, and the epoch-millis
# timestamp), so force those specific cells to text first to keep them
# stored verbatim instead of being converted to a number/date serial.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "231249"

$ws.Range("B2").Value = "Histology Lab"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "03/05/2025"

$ws.Range("D2").Value = "01:25:37"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1746224737580"

# Column F (old "Type" column) is no longer part of the table - remove it
$ws.Range("F1:F2").Clear()
